$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; all existing rows from 41 downward shift
# down by one (old row 41 becomes row 42, ..., old row 115 becomes row 116).
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point.
# Columns A-C and E-M carry the same values as the row that was pushed down
# (now row 42); only D (fecha) and N-T (price/unit/origin columns) differ.
$ws.Cells.Item(41, 1).Value = 8
$ws.Cells.Item(41, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41, 3).Value = "Coquimbo"
$ws.Cells.Item(41, 4).Value = 44665
$ws.Cells.Item(41, 5).Value = 4
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100109
$ws.Cells.Item(41, 8).Value = "Uva"
$ws.Cells.Item(41, 9).Value = 100109001
$ws.Cells.Item(41, 10).Value = "Uva"
$ws.Cells.Item(41, 11).Value = "Red Globe"
$ws.Cells.Item(41, 12).Value = "Primera"
$ws.Cells.Item(41, 13).Value = 400
$ws.Cells.Item(41, 14).Value = 9000
$ws.Cells.Item(41, 15).Value = 10000
$ws.Cells.Item(41, 16).Value = 9500
$ws.Cells.Item(41, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(41, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(41, 19).Value = 528
$ws.Cells.Item(41, 20).Value = 18
